$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.406.67"
$ws.Range("E2").Value = "  -0.68%  "
$ws.Range("D3").Value = "3.077.44"
$ws.Range("E3").Value = "  -1.30%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "574.91"
$ws.Range("E5").Value = "  -0.43%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "171.30"
$ws.Range("E6").Value = "  -0.52%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").Value = "3.075.13"
$ws.Range("E8").Value = "  -1.29%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.510"
$ws.Range("E9").Value = "  -2.03%  "
$ws.Range("E10").Value = "  -1.84%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.150"
$ws.Range("E11").Value = "  -2.16%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.470"
$ws.Range("E12").Value = "  -2.50%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000239"
$ws.Range("E13").Value = "  -3.74%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.77"
$ws.Range("E14").Value = "  -3.72%  "
$ws.Range("E15").Value = "  -1.47%  "
$ws.Range("D16").Value = "3.589.90"
$ws.Range("E16").Value = "  -1.22%  "
$ws.Range("D17").Value = "66.330.20"
$ws.Range("E17").Value = "  -0.77%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.96"
$ws.Range("E18").Value = "  -2.51%  "
$ws.Range("B19").Value = "WrappedEther"
$ws.Range("C19").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D19").Value = "3.076.82"
$ws.Range("E19").Value = "  -1.26%  "
$ws.Range("B20").Value = "Chainlink"
$ws.Range("C20").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "16.55"
$ws.Range("E20").Value = "  +2.02%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "485.55"
$ws.Range("E21").Value = "  +2.40%  "
$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.69"
$ws.Range("E22").Value = "  -2.43%  "
$ws.Range("B23").Value = "Polygon"
$ws.Range("C23").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.686"
$ws.Range("E23").Value = "  -3.19%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "82.33"
$ws.Range("E24").Value = "  -1.62%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.66"
$ws.Range("E25").Value = "  -4.33%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.22"
$ws.Range("E26").Value = "  -2.55%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.13"
$ws.Range("E27").Value = "  -1.36%  "
$ws.Range("E28").Value = "  +0.03%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.93"
$ws.Range("E29").Value = "  +0.20%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.26"
$ws.Range("E30").Value = "  -4.85%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.60"
$ws.Range("E31").Value = "  -2.85%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "27.77"
$ws.Range("E32").Value = "  -2.65%  "
$ws.Range("E33").Value = "  -2.89%  "
$ws.Range("D34").Value = "0.0₃0910"
$ws.Range("E34").Value = "  -4.11%  "
$ws.Range("E35").Value = "  -0.08%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "48.28"
$ws.Range("E36").Value = "  +2.70%  "
$ws.Range("B37").Value = "Filecoin"
$ws.Range("C37").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.56"
$ws.Range("E37").Value = "  -4.68%  "
$ws.Range("B38").Value = "Mantle"
$ws.Range("C38").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.940"
$ws.Range("E38").Value = "  -3.52%  "
$ws.Range("E39").Value = "  -1.28%  "
$ws.Range("B40").Value = "TheGraph"
$ws.Range("C40").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.302"
$ws.Range("E40").Value = "  -3.01%  "
$ws.Range("B41").Value = "Stacks"
$ws.Range("C41").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.96"
$ws.Range("E41").Value = "  -4.58%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.23"
$ws.Range("E42").Value = "  -4.39%  "
$ws.Range("D43").Value = "2.769.51"
$ws.Range("E43").Value = "  -1.52%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.52"
$ws.Range("E44").Value = "  -0.76%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0343"
$ws.Range("E45").Value = "  -2.72%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "134.64"
$ws.Range("E46").Value = "  -0.63%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "365.12"
$ws.Range("E47").Value = "  -4.48%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "24.32"
$ws.Range("E49").Value = "  -2.21%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.15"
$ws.Range("E50").Value = "  -2.05%  "
$ws.Range("E51").Value = "  -1.99%  "
